$wb = $excel.ActiveWorkbook

# --- 1. Insert a new "State" column into hotel_info, right after "Hotel_Name" ---
$hotel = $wb.Worksheets.Item("hotel_info")

# Insert a new column before column C (City), which becomes the new "State" column (column C)
$hotel.Columns.Item(3).Insert()

$hotel.Range("C1").Value = "State"
$hotel.Range("C2").Value = "Louisiana"

# --- 2. Reorder worksheets so review_info comes before hotel_info ---
$review = $wb.Worksheets.Item("review_info")
$review.Move($hotel)

$wb.Save()
